$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.531.75"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.825.43"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.61"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5104"
$ws.Range("E7").Value = "  -5.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08198"
$ws.Range("E9").Value = "  +6.20%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.361"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.564"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "1.823.50"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.84"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06666"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.84"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.113"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "28.572.65"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.18"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "2.030.24"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.417"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.53"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1086"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.770"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.664"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07052"
$ws.Range("E35").Value = "  -4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2232"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02358"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.795"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6374"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.29"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.401"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.57"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5981"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.30"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.993"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.197"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06943"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +3.93%  "
